$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user registration row appended by the login/registration validation feature.
$ws.Cells.Item(4, 1).Value = "Sebastián"
$ws.Cells.Item(4, 2).Value = "Palacio"
$ws.Cells.Item(4, 3).Value = 1000762620
$ws.Cells.Item(4, 4).Value = "sebasx200"

$pwd = $ws.Cells.Item(4, 5)
$pwd.NumberFormat = "@"
$pwd.Value = "1234"
$pwd.Style = "Normal"

$ws.Cells.Item(4, 6).Value = "sebastian_palacio23231@elpoli,edu,co"
$ws.Cells.Item(4, 7).Value = "No tiene"
